# Generate Report for Handback
# Update the handback status report for 86d8f3d0-c847-43ef-9a15-69810fdd18af:
#  - Status flips from "ht" (human translation) to "mt" (machine translation)
#  - Handoff / Handback timestamps are refreshed to the new generation run

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 86d8f3d0 file (rows 3 and 5 both
# reference the same generation run for that source file group)
$wsOverview.Range("G3").Value = "2016-08-27 10:17:17"
$wsOverview.Range("G5").Value = "2016-08-27 10:17:17"

# zh-cn sheet: Status + Correspond Handoff/Handback datetimes for the 86d8f3d0 file
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-27 10:17:12"
$wsZhCn.Range("H5").Value = "2016-08-27 10:17:12"
$wsZhCn.Range("K3").Value = "2016-08-27 10:17:29"
$wsZhCn.Range("K5").Value = "2016-08-27 10:17:29"

# de-de sheet: Status + Correspond Handback datetime for the 86d8f3d0 file
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-27 10:17:17"
$wsDeDe.Range("H5").Value = "2016-08-27 10:17:17"
$wsDeDe.Range("K3").Value = "2016-08-27 10:17:36"
$wsDeDe.Range("K5").Value = "2016-08-27 10:17:36"
